$d = $word.ActiveDocument

# Locate the paragraph containing the old comment line.
$targetIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*# Password validation logic*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph '# Password validation logic'"
}

$para = $d.Paragraphs($targetIndex)

# Replace its text with the new first line of the code block.
$para.Range.Text = "    # Function to validate password "

# Remaining lines of the new code block, to be inserted as new
# paragraphs immediately following the updated paragraph, in order.
$newLines = @(
    "    # per given requirements",
    "    if len(password) < 8:",
    "        return False",
    "    if not any(char.isdigit() for char in password):",
    "        return False",
    "    if not any(char in ['!', '@', '#', '$', '%', '^', '&', '*'] for char in password):",
    "        return False"
)

$insertAfterIndex = $targetIndex
foreach ($line in $newLines) {
    $d.Paragraphs($insertAfterIndex).Range.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    $d.Paragraphs($insertAfterIndex).Range.Text = $line
}

Write-Output "Inserted $($newLines.Count) new paragraphs after paragraph $targetIndex"
